# Insert a new "raw accuracy" column before the existing results columns
# (old column C "BAARD2" ... shifts to D, etc.) and populate the new
# column C with per-attack raw accuracy values plus a header label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting column C pushes C:H -> D:I and copies formatting from column B,
# which already matches the bordered/bold/centered style used throughout
# column A/B/C in this sheet.
$ws.Columns.Item(3).Insert()

# The old (empty) B1/B2 placeholder cells are not reused for the new column;
# the new C1/C2 cells start out completely blank instead, so drop whatever
# leftover content/formatting the insert left behind in B1/B2.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# New header cell for the inserted column.
$ws.Range("C3").Value = "Accuracy after attack"

# Raw accuracy values for each attack/epsilon row.
$ws.Range("C4").Value  = 80.80000000000001
$ws.Range("C5").Value  = 0
$ws.Range("C6").Value  = 0
$ws.Range("C7").Value  = 0
$ws.Range("C8").Value  = 100
$ws.Range("C9").Value  = 99.7
$ws.Range("C10").Value = 99.3
$ws.Range("C11").Value = 99.09999999999999
$ws.Range("C12").Value = 98.40000000000001
$ws.Range("C13").Value = 96.40000000000001
$ws.Range("C14").Value = 99.7
$ws.Range("C15").Value = 97.59999999999999
$ws.Range("C16").Value = 91.59999999999999
$ws.Range("C17").Value = 69.39999999999999
$ws.Range("C18").Value = 30.3
$ws.Range("C19").Value = 10.2
$ws.Range("C20").Value = 2.3
$ws.Range("C21").Value = 0.5
$ws.Range("C22").Value = 25
$ws.Range("C23").Value = 3.5
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 33.1
$ws.Range("C27").Value = 60
$ws.Range("C28").Value = 1.6
$ws.Range("C29").Value = 84.2
$ws.Range("C30").Value = 6.9
$ws.Range("C31").Value = 0.8
$ws.Range("C32").Value = 0.4
$ws.Range("C33").Value = 0.4
$ws.Range("C34").Value = 99.59999999999999
$ws.Range("C35").Value = 99
$ws.Range("C36").Value = 96.8
$ws.Range("C37").Value = 99.8
$ws.Range("C38").Value = 99.8
$ws.Range("C39").Value = 98.40000000000001
$ws.Range("C40").Value = 93.59999999999999
$ws.Range("C41").Value = 74
